$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while forcing text interpretation (preserves original
# number-less default style, avoiding Excel auto-converting numeric-looking
# strings like "1.004" or "0.00000000117" into actual numbers).
function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = $origStyle
}

# Row 2
Set-TextValue "D2" "29.298.43"
$ws.Range("E2").Value = "  -0.27%  "

# Row 3
Set-TextValue "D3" "1.831.34"
$ws.Range("E3").Value = "  -0.56%  "

# Row 4
Set-TextValue "D4" "1.004"
$ws.Range("E4").Value = "  +0.49%  "

# Row 5
Set-TextValue "D5" "235.56"
$ws.Range("E5").Value = "  -1.56%  "

# Row 6
Set-TextValue "D6" "0.6034"
$ws.Range("E6").Value = "  -3.59%  "

# Row 7
Set-TextValue "D7" "1.005"
$ws.Range("E7").Value = "  +0.40%  "

# Row 8
Set-TextValue "D8" "0.07050"
$ws.Range("E8").Value = "  -5.06%  "

# Row 9
Set-TextValue "D9" "0.2789"
$ws.Range("E9").Value = "  -3.60%  "

# Row 10
Set-TextValue "D10" "23.52"
$ws.Range("E10").Value = "  -5.62%  "

# Row 11
Set-TextValue "D11" "0.07668"
$ws.Range("E11").Value = "  -0.66%  "

# Row 12
Set-TextValue "D12" "1.831.57"
$ws.Range("E12").Value = "  -0.81%  "

# Row 13
Set-TextValue "D13" "4.795"
$ws.Range("E13").Value = "  -3.57%  "

# Row 14
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue "D14" "0.6286"
$ws.Range("E14").Value = "  -6.74%  "

# Row 15
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue "D15" "0.000009925"
$ws.Range("E15").Value = "  -3.19%  "

# Row 16
Set-TextValue "D16" "79.15"
$ws.Range("E16").Value = "  -3.22%  "

# Row 17
Set-TextValue "D17" "29.299.85"
$ws.Range("E17").Value = "  -0.28%  "

# Row 18
Set-TextValue "D18" "5.846"
$ws.Range("E18").Value = "  -5.80%  "

# Row 19
Set-TextValue "D19" "224.42"
$ws.Range("E19").Value = "  -4.02%  "

# Row 20
Set-TextValue "D20" "1.005"
$ws.Range("E20").Value = "  +0.44%  "

# Row 21
Set-TextValue "D21" "11.71"
$ws.Range("E21").Value = "  -4.88%  "

# Row 22
Set-TextValue "D22" "7.015"
$ws.Range("E22").Value = "  -3.54%  "

# Row 23
Set-TextValue "D23" "1.006"
$ws.Range("E23").Value = "  +0.45%  "

# Row 24
Set-TextValue "D24" "155.62"
$ws.Range("E24").Value = "  -1.69%  "

# Row 25
Set-TextValue "D25" "7.983"
$ws.Range("E25").Value = "  -5.82%  "

# Row 26
Set-TextValue "D26" "0.1299"
$ws.Range("E26").Value = "  -3.49%  "

# Row 27
Set-TextValue "D27" "16.57"
$ws.Range("E27").Value = "  -4.22%  "

# Row 28
Set-TextValue "D28" "1.478"
$ws.Range("E28").Value = "  +1.00%  "

# Row 29
Set-TextValue "D29" "0.06337"
$ws.Range("E29").Value = "  -12.19%  "

# Row 30
Set-TextValue "D30" "1.448"
$ws.Range("E30").Value = "  -2.06%  "

# Row 31
Set-TextValue "D31" "3.847"
$ws.Range("E31").Value = "  -4.37%  "

# Row 32
Set-TextValue "D32" "3.804"
$ws.Range("E32").Value = "  -6.03%  "

# Row 33
Set-TextValue "D33" "1.109"
$ws.Range("E33").Value = "  -2.64%  "

# Row 34
Set-TextValue "D34" "1.735"
$ws.Range("E34").Value = "  -4.56%  "

# Row 35
Set-TextValue "D35" "0.6460"
$ws.Range("E35").Value = "  -7.32%  "

# Row 36
Set-TextValue "D36" "2.550"
$ws.Range("E36").Value = "  -0.85%  "

# Row 37
Set-TextValue "D37" "1.220.43"
$ws.Range("E37").Value = "  -0.90%  "

# Row 38
Set-TextValue "D38" "2.746"
$ws.Range("E38").Value = "  -2.49%  "

# Row 39
Set-TextValue "D39" "0.01738"
$ws.Range("E39").Value = "  -5.59%  "

# Row 40
Set-TextValue "D40" "6.504"
$ws.Range("E40").Value = "  -5.96%  "

# Row 41
Set-TextValue "D41" "0.8985"
$ws.Range("E41").Value = "  -6.90%  "

# Row 42
$ws.Range("E42").Value = "  +0.42%  "

# Row 43
Set-TextValue "D43" "1.991.53"
$ws.Range("E43").Value = "  -0.21%  "

# Row 44
Set-TextValue "D44" "100.52"
$ws.Range("E44").Value = "  -0.47%  "

# Row 45
Set-TextValue "D45" "62.61"
$ws.Range("E45").Value = "  -4.32%  "

# Row 46
Set-TextValue "D46" "0.00000000117"
$ws.Range("E46").Value = "  -2.20%  "

# Row 47
Set-TextValue "D47" "8.538"
$ws.Range("E47").Value = "  -4.07%  "

# Row 48
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D48" "1.579"
$ws.Range("E48").Value = "  -8.03%  "

# Row 49
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue "D49" "0.4563"
$ws.Range("E49").Value = "  -0.38%  "

# Row 50
Set-TextValue "D50" "0.05500"
$ws.Range("E50").Value = "  -2.80%  "

# Row 51
Set-TextValue "D51" "6.403"
$ws.Range("E51").Value = "  -7.82%  "
